# Add 2022-Q3 data:
#  - Insert a new worksheet "2022-Q3" (copied from "2022-Q2" so it keeps the
#    same layout/formatting) positioned right after "总计" and before
#    "2022-Q2".
#  - Fill it with the 2022-Q3 fund-holdings table (25 data rows).
#  - Update the "总计" (summary) sheet: insert a new top data row for
#    2022-Q3 and shift the existing quarters down by one row.

$wb = $excel.ActiveWorkbook

function Fill-DataRow {
    param($ws, [int]$r, $b, $c, $d, $e, $f, $g, $h)
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
    $ws.Cells.Item($r, 6).Value = $f
    $ws.Cells.Item($r, 7).Value = $g
    $ws.Cells.Item($r, 8).Value = $h
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet by duplicating "2022-Q2" (sheet #2),
#    so it keeps identical headers/column-widths/cell styling. The copy
#    is placed immediately before "2022-Q2" -> position 2, right after
#    "总计".
# ---------------------------------------------------------------------
$q2Sheet = $wb.Worksheets.Item(2)
$q2Sheet.Copy($q2Sheet)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "2022-Q3"

# The source sheet ("2022-Q2") only has 16 data rows (rows 2-17); the
# 2022-Q3 table needs 25 data rows (rows 2-26), so extend column A's
# formatting (bold/border/center style) down to row 26.
$newSheet.Range("A17").Copy()
$newSheet.Range("A18:A26").PasteSpecial(-4122)
for ($r = 18; $r -le 26; $r++) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2
}

# Columns B-G hold text values (fund code / name / size / position %'s)
# in this workbook, even though most look numeric - force text format
# before assigning so they are stored as text, not re-interpreted as
# numbers. (Row 25/26 column G are a genuine numeric 0, handled below.)
$newSheet.Range("B2:G24").NumberFormat = "@"
$newSheet.Range("B25:F26").NumberFormat = "@"

Fill-DataRow $newSheet 2  "515210" "国泰中证钢铁ETF" "14.23" "97.88" "9.19" "1.3077" 2
Fill-DataRow $newSheet 3  "502023" "鹏华国证钢铁行业指数（LOF）A" "9.48" "94.49" "12.84" "1.2172" 1
Fill-DataRow $newSheet 4  "012810" "鹏华国证钢铁行业指数（LOF）C" "4.34" "94.49" "12.84" "0.5573" 1
Fill-DataRow $newSheet 5  "168203" "中融国证钢铁行业指数A" "3.34" "92.81" "12.62" "0.4215" 1
Fill-DataRow $newSheet 6  "501029" "华宝标普中国A股红利机会指数（LOF）A" "8.11" "94.26" "1.39" "0.1127" 9
Fill-DataRow $newSheet 7  "393001" "中海优势精选灵活配置混合" "1.50" "78.67" "6.81" "0.1022" 8
Fill-DataRow $newSheet 8  "013934" "长江红利回报混合A" "2.39" "87.78" "4.04" "0.0966" 4
Fill-DataRow $newSheet 9  "510160" "南方中证南方小康产业ETF" "2.34" "99.43" "2.83" "0.0662" 4
Fill-DataRow $newSheet 10 "005125" "华宝标普中国A股红利机会指数C" "3.38" "94.26" "1.39" "0.0470" 9
Fill-DataRow $newSheet 11 "003845" "汇安丰恒灵活配置混合A" "1.27" "61.83" "2.93" "0.0372" 8
Fill-DataRow $newSheet 12 "510760" "国泰上证综合ETF" "2.75" "94.96" "1.18" "0.0324" 10
Fill-DataRow $newSheet 13 "517180" "南方富时中国国企开放共赢ETF" "0.36" "97.73" "3.83" "0.0138" 6
Fill-DataRow $newSheet 14 "515500" "海富通中证长三角领先ETF" "0.30" "97.12" "4.36" "0.0131" 8
Fill-DataRow $newSheet 15 "660006" "农银大盘蓝筹混合" "1.27" "84.42" "1.02" "0.0130" 10
Fill-DataRow $newSheet 16 "159719" "平安富时中国国企开放共赢ETF" "0.34" "94.14" "3.76" "0.0128" 6
Fill-DataRow $newSheet 17 "013935" "长江红利回报混合C" "0.24" "87.78" "4.04" "0.0097" 4
Fill-DataRow $newSheet 18 "013802" "财通资管中证钢铁指数A" "0.08" "92.45" "6.19" "0.0050" 5
Fill-DataRow $newSheet 19 "159944" "广发中证全指原材料ETF" "0.20" "98.02" "1.13" "0.0023" 9
Fill-DataRow $newSheet 20 "004403" "平安股息精选沪港深股票A" "0.08" "91.74" "2.52" "0.0020" 9
Fill-DataRow $newSheet 21 "006347" "安信量化优选股票C" "0.15" "90.50" "0.81" "0.0012" 6
Fill-DataRow $newSheet 22 "013803" "财通资管中证钢铁指数C" "0.01" "92.45" "6.19" "0.0006" 5
Fill-DataRow $newSheet 23 "004404" "平安股息精选沪港深股票C" "0.02" "91.74" "2.52" "0.0005" 9
Fill-DataRow $newSheet 24 "006346" "安信量化优选股票A" "0.03" "90.50" "0.81" "0.0002" 6
Fill-DataRow $newSheet 25 "003846" "汇安丰恒灵活配置混合C" "0.00" "61.83" "2.93" 0 8
Fill-DataRow $newSheet 26 "016815" "中融国证钢铁行业指数C" "0.00" "92.81" "12.62" 0 1

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet: shift the existing 7 quarters down
#    one row and write the new 2022-Q3 totals into row 2.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Row 9 is brand new - copy row 8's column-A formatting down to it first.
$summary.Range("A8").Copy()
$summary.Range("A9").PasteSpecial(-4122)

$summary.Cells.Item(9, 1).Value = 7
$summary.Cells.Item(9, 2).Value = "2020-Q4"
$summary.Cells.Item(9, 3).Value = 24
$summary.Cells.Item(9, 4).Value = 4.8

$summary.Cells.Item(8, 1).Value = 6
$summary.Cells.Item(8, 2).Value = "2021-Q1"
$summary.Cells.Item(8, 3).Value = 73
$summary.Cells.Item(8, 4).Value = 22.38

$summary.Cells.Item(7, 1).Value = 5
$summary.Cells.Item(7, 2).Value = "2021-Q2"
$summary.Cells.Item(7, 3).Value = 27
$summary.Cells.Item(7, 4).Value = 15.05

$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(6, 2).Value = "2021-Q3"
$summary.Cells.Item(6, 3).Value = 53
$summary.Cells.Item(6, 4).Value = 24.14

$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(5, 2).Value = "2021-Q4"
$summary.Cells.Item(5, 3).Value = 44
$summary.Cells.Item(5, 4).Value = 15.42

$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(4, 2).Value = "2022-Q1"
$summary.Cells.Item(4, 3).Value = 38
$summary.Cells.Item(4, 4).Value = 12.27

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = "2022-Q2"
$summary.Cells.Item(3, 3).Value = 16
$summary.Cells.Item(3, 4).Value = 5

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 25
$summary.Cells.Item(2, 4).Value = 4.07
